$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.264055609703064
$ws.Range("B1").Value = 2.74508810043335
$ws.Range("C1").Value = 5.105611801147461
$ws.Range("D1").Value = 2.005113363265991
$ws.Range("E1").Value = 1.033000349998474
